$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'98.507.24"
$ws.Range("E2").Value = "'  +0.07%  "
$ws.Range("D3").Value = "'3.354.33"
$ws.Range("E3").Value = "'  +0.78%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'257.45"
$ws.Range("E5").Value = "'  -0.18%  "
$ws.Range("D6").Value = "'663.59"
$ws.Range("E6").Value = "'  +6.12%  "
$ws.Range("E7").Value = "'  +9.28%  "
$ws.Range("D8").Value = "'0.477"
$ws.Range("E8").Value = "'  +22.98%  "
$ws.Range("E9").Value = "'  +24.45%  "
$ws.Range("E10").Value = "'  -0.02%  "
$ws.Range("D11").Value = "'3.352.85"
$ws.Range("E11").Value = "'  +0.76%  "
$ws.Range("D12").Value = "'0.215"
$ws.Range("E12").Value = "'  +8.37%  "
$ws.Range("D13").Value = "'42.16"
$ws.Range("E13").Value = "'  +13.76%  "
$ws.Range("D14").Value = "'0.0000275"
$ws.Range("E14").Value = "'  +11.17%  "
$ws.Range("D15").Value = "'98.623.53"
$ws.Range("E15").Value = "'  +0.47%  "
$ws.Range("E16").Value = "'  +3.72%  "
$ws.Range("D17").Value = "'3.976.86"
$ws.Range("E17").Value = "'  +0.98%  "
$ws.Range("D18").Value = "'3.345.69"
$ws.Range("E18").Value = "'  +0.69%  "
$ws.Range("D19").Value = "'7.62"
$ws.Range("E19").Value = "'  +26.06%  "
$ws.Range("D20").Value = "'16.71"
$ws.Range("E20").Value = "'  +10.83%  "
$ws.Range("D21").Value = "'3.58"
$ws.Range("E21").Value = "'  +1.48%  "
$ws.Range("D22").Value = "'528.89"
$ws.Range("E22").Value = "'  +8.86%  "
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "'  +13.64%  "
$ws.Range("E24").Value = "'  +4.33%  "
$ws.Range("D25").Value = "'0.434"
$ws.Range("E25").Value = "'  +52.14%  "
$ws.Range("D26").Value = "'102.22"
$ws.Range("E26").Value = "'  +15.46%  "
$ws.Range("D27").Value = "'6.09"
$ws.Range("E27").Value = "'  +8.62%  "
$ws.Range("D28").Value = "'12.52"
$ws.Range("E28").Value = "'  +6.09%  "
$ws.Range("D29").Value = "'3.533.28"
$ws.Range("E29").Value = "'  +0.86%  "
$ws.Range("E30").Value = "'  +6.56%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "'  -0.02%  "
$ws.Range("D32").Value = "'11.03"
$ws.Range("E32").Value = "'  +15.05%  "
$ws.Range("E33").Value = "'  -2.19%  "
$ws.Range("E34").Value = "'  -0.16%  "
$ws.Range("D35").Value = "'29.33"
$ws.Range("E35").Value = "'  +5.78%  "
$ws.Range("D36").Value = "'0.539"
$ws.Range("E36").Value = "'  +17.69%  "
$ws.Range("D37").Value = "'7.80"
$ws.Range("E37").Value = "'  +7.30%  "
$ws.Range("E38").Value = "'  +8.96%  "
$ws.Range("E39").Value = "'  +5.69%  "
$ws.Range("D40").Value = "'525.37"
$ws.Range("E40").Value = "'  +6.38%  "
$ws.Range("B41").Value = "'Fetch.AI"
$ws.Range("C41").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'1.33"
$ws.Range("E41").Value = "'  +6.20%  "
$ws.Range("B42").Value = "'WhiteBITCoin"
$ws.Range("C42").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'24.71"
$ws.Range("E42").Value = "'  -0.53%  "
$ws.Range("D43").Value = "'3.87"
$ws.Range("E43").Value = "'  +4.68%  "
$ws.Range("E44").Value = "'  +32.74%  "
$ws.Range("D45").Value = "'3.43"
$ws.Range("E45").Value = "'  +4.51%  "
$ws.Range("D46").Value = "'0.820"
$ws.Range("E46").Value = "'  +5.55%  "
$ws.Range("E47").Value = "'  +0.04%  "
$ws.Range("B48").Value = "'Stacks"
$ws.Range("C48").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.06"
$ws.Range("E48").Value = "'  +7.16%  "
$ws.Range("B49").Value = "'Filecoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'5.12"
$ws.Range("E49").Value = "'  +11.50%  "
$ws.Range("D50").Value = "'7.81"
$ws.Range("E50").Value = "'  +16.70%  "
$ws.Range("B51").Value = "'ImmutableX"
$ws.Range("C51").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "'1.53"
$ws.Range("E51").Value = "'  +12.58%  "
